$d = $word.ActiveDocument

# 1. Update the letter date.
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 2)

# 2. Split the mailing address line into two paragraphs (street / city-state-zip)
#    with a blank paragraph after, matching the new mailing-address block.
$d.Content.Find.Execute("521 Rough  Ready Road, San Jose CA 95133", $true, $false, $false, $false, $false,
                         $true, 1, $false, "521 Rough  Ready Road^pSan Jose, CA 95133^p", 2)

# 3. Remove the two blank paragraphs that used to sit right after
#    "... Board of Directors", leaving just the single trailing blank
#    "Title" styled paragraph that was already further down.
$find = $d.Content.Find
$find.Execute("Board of Directors") | Out-Null
$rng = $find.Parent
$rng.Collapse(0)
[void]$rng.MoveStart(1, 1)
$target = $rng.Start

$count = $d.Paragraphs.Count
$idx = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Start -eq $target) {
        $idx = $i
        break
    }
}
if ($idx -ne -1) {
    $d.Paragraphs.Item($idx).Range.Delete()
    $d.Paragraphs.Item($idx).Range.Delete()
}
